# Reflects Eddie Bernice Johnson (D, TX-30) moving from the "NO" / holdout
# list to supporting impeachment, updating the dependent summary tables
# ("8 holdouts" instead of 9).

$wb = $excel.ActiveWorkbook

# --- prezresults2016: drop the D/NO=1 row, bump D/YES from 203 to 204 ---
$ws1 = $wb.Worksheets.Item("prezresults2016")
$ws1.Rows.Item(2).Delete()
$ws1.Range("C3").Value = 204

# --- full_list_of_nos: remove Eddie Bernice Johnson's row (was row 6) ---
$ws12 = $wb.Worksheets.Item("full_list_of_nos")
$ws12.Rows.Item(6).Delete()

# --- gdp_vs_nationalavg: NO/ABOVE 2->1, YES/ABOVE 131->132 ---
$ws2 = $wb.Worksheets.Item("gdp_vs_nationalavg")
$ws2.Range("C2").Value = 1
$ws2.Range("C4").Value = 132

# --- college_vs_nationalavg: NO/BELOW 8->7, YES/BELOW 92->93 ---
$ws3 = $wb.Worksheets.Item("college_vs_nationalavg")
$ws3.Range("C3").Value = 7
$ws3.Range("C5").Value = 93

# --- nonwhite_vs_nationalavg: NO/ABOVE 3->2, YES/ABOVE 141->142 ---
$ws4 = $wb.Worksheets.Item("nonwhite_vs_nationalavg")
$ws4.Range("C2").Value = 2
$ws4.Range("C4").Value = 142

# --- rural_morethanfifth: NO/BELOW 4->3, YES/BELOW 191->192 ---
$ws5 = $wb.Worksheets.Item("rural_morethanfifth")
$ws5.Range("C3").Value = 3
$ws5.Range("C5").Value = 192

# --- margin_5_or_less: NO/more_than_5_points 3->2, YES/more_than_5_points 206->207 ---
$ws6 = $wb.Worksheets.Item("margin_5_or_less")
$ws6.Range("C3").Value = 2
$ws6.Range("C5").Value = 207
